$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Remove the four obsolete paragraphs that used to sit right after
# "For publishers admin should be able to add bids manually ..." and its
# trailing blank paragraph:
#   "min starting budjet 50$ and min cpm or cpc 0.001$"
#   "Modify targeting"
#   "Admin can also run in-house campaign like advitisers"
#   "Multiple cateogry. Cateogry while adding site"
# Deleting a paragraph's Range (which includes its end-of-paragraph mark)
# removes the whole paragraph and renumbers everything after it, so deleting
# the same index four times in a row removes all four.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(42).Range.Delete()
$d.Paragraphs.Item(42).Range.Delete()
$d.Paragraphs.Item(42).Range.Delete()
$d.Paragraphs.Item(42).Range.Delete()

# Paragraph 42 is now "Reports for both advitisers and publishers".

# ---------------------------------------------------------------------------
# Add the new "Country Targeting / and domain targeting / does not work."
# paragraph right after it. It is built as three short paragraphs (one per
# run) which are then joined back into a single paragraph by deleting the
# paragraph marks between them -- that merge keeps the three pieces of text
# as three distinct runs instead of Word coalescing them into one.
# ---------------------------------------------------------------------------
$d.Paragraphs.Item(42).Range.InsertParagraphAfter()
$d.Paragraphs.Item(43).Range.InsertAfter("Country Targeting ")
$d.Paragraphs.Item(43).Range.InsertParagraphAfter()
$d.Paragraphs.Item(44).Range.InsertAfter("and domain targeting ")
$d.Paragraphs.Item(44).Range.InsertParagraphAfter()
$d.Paragraphs.Item(45).Range.InsertAfter("does not work. ")

$mark = $d.Paragraphs.Item(43).Range.End
$d.Range($mark - 1, $mark).Delete()
$mark = $d.Paragraphs.Item(43).Range.End
$d.Range($mark - 1, $mark).Delete()

# Paragraph 43 is now the merged "Country Targeting and domain targeting
# does not work." paragraph (three runs). Paragraph 44 is the pre-existing
# blank, right-aligned paragraph.

# ---------------------------------------------------------------------------
# Five more blank, right-aligned paragraphs are needed after it.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt 5; $i++) {
    $d.Paragraphs.Item(44).Range.InsertParagraphAfter()
    $d.Paragraphs.Item(45).Format.Alignment = 2
}

# ---------------------------------------------------------------------------
# Paragraph 50 (currently blank) becomes the new "Admin can also run
# in-house campaign like advitisers -> For this create a publisher account
# and create manual deposit and run it." paragraph, made of two runs. A
# scratch paragraph is inserted after it to hold the second run's text, and
# then the two paragraph marks are merged back together -- leaving the
# blank paragraph that originally followed untouched.
# ---------------------------------------------------------------------------
$arrow = [string][char]0x2192
$d.Paragraphs.Item(50).Range.InsertAfter("Admin can also run in-house campaign like advitisers " + $arrow)
$d.Paragraphs.Item(50).Range.InsertParagraphAfter()
$d.Paragraphs.Item(51).Range.InsertAfter(" For this create a publisher account and create manual deposit and run it.")

$mark = $d.Paragraphs.Item(50).Range.End
$d.Range($mark - 1, $mark).Delete()
